$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 8 with the new "Red Coat" entry
$ws.Range("C8").Value = "Red Coat"
$ws.Range("D8").Value = "https://movieleatherjackets.com/demi-lovato-i-love-me-song-red-coat/"
$ws.Range("E8").Value = "Red Coat"
$ws.Range("F8").Value = "demi_lovato_red_coat.png"

$ws.Hyperlinks.Add($ws.Range("D8"), "https://movieleatherjackets.com/demi-lovato-i-love-me-song-red-coat/") | Out-Null
$ws.Range("D8").Style = $ws.Range("D2").Style

# Move the active selection to D14 as in the edited file
$ws.Range("D14").Select()
